$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '61.111.25'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +0.66%  '
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.928.59'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +0.76%  '
$c.Style = 'Normal'
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '594.01'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  +1.24%  '
$c.Style = 'Normal'
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '146.22'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -1.22%  '
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +3.03%  '
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +0.52%  '
$c.Style = 'Normal'
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.441'
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -1.80%  '
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +1.16%  '
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '33.83'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -0.84%  '
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.Style = 'Normal'
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '3.413.22'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  +0.72%  '
$c.Style = 'Normal'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '61.088.52'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  +0.66%  '
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '6.72'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -1.45%  '
$c.Style = 'Normal'
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '2.930.95'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.79%  '
$c.Style = 'Normal'
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '432.06'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  +1.10%  '
$c.Style = 'Normal'
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '13.47'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -1.27%  '
$c.Style = 'Normal'
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '0.683'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +1.91%  '
$c.Style = 'Normal'
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '7.10'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  -0.34%  '
$c.Style = 'Normal'
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '81.39'
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c.Style = 'Normal'
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '11.06'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.25%  '
$c.Style = 'Normal'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  +0.33%  '
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +1.75%  '
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.Style = 'Normal'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +6.19%  '
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.10%  '
$c.Style = 'Normal'
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '7.10'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -1.98%  '
$c.Style = 'Normal'
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '26.48'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +1.31%  '
$c.Style = 'Normal'
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '0.0₃0857'
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +2.01%  '
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +0.21%  '
$c.Style = 'Normal'
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '5.64'
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.Style = 'Normal'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '3.07'
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  +3.64%  '
$c.Style = 'Normal'
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.124'
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +2.32%  '
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -1.71%  '
$c.Style = 'Normal'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '8.60'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -1.51%  '
$c.Style = 'Normal'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.288'
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  -1.48%  '
$c.Style = 'Normal'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '39.72'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -4.38%  '
$c.Style = 'Normal'
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '376.24'
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +1.49%  '
$c.Style = 'Normal'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.0346'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '2.720.04'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +2.44%  '
$c.Style = 'Normal'
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '130.79'
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -1.96%  '
$c.Style = 'Normal'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '24.17'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -3.97%  '
$c.Style = 'Normal'
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '2.03'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -2.87%  '
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  +2.66%  '
$c.Style = 'Normal'
